# This workbook has a single worksheet "Avverkningsanmälningar".
# The commit bumps the "Förändrad" column (C) from 2023-09-20 (45189)
# to 2023-09-21 (45190) for every existing data row (2..233), and
# appends one new data row (234) for case "A 44353-2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump column C ("Förändrad") for all existing data rows to 45190.
$ws.Range("C2:C233").Value = 45190

# 2) The last existing row (233) gains an explicit row height in the
#    saved XML (ht="15" customHeight="1"), matching the other rows.
$ws.Rows.Item(233).RowHeight = 15

# 3) Append the new row (234) with the new case's data.
$ws.Range("A234").Value = "A 44353-2023"
$ws.Range("B234").Value = 45188
$ws.Range("C234").Value = 45190
$ws.Range("D234").Value = "HALLANDS LÄN"
$ws.Range("E234").Value = "KUNGSBACKA"
$ws.Range("G234").Value = 4.8
$ws.Range("H234").Value = 0
$ws.Range("I234").Value = 0
$ws.Range("J234").Value = 0
$ws.Range("K234").Value = 0
$ws.Range("L234").Value = 0
$ws.Range("M234").Value = 0
$ws.Range("N234").Value = 0
$ws.Range("O234").Value = 0
$ws.Range("P234").Value = 0
$ws.Range("Q234").Value = 0

# 4) Match formatting of the new row to the row above it: date style
#    for B/C, and the wrap-text style used (empty) in column R.
$ws.Range("B233:C233").Copy()
$ws.Range("B234:C234").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("R233").Copy()
$ws.Range("R234").PasteSpecial(-4122)  # xlPasteFormats
